$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.075827121734619
$ws.Range("B1").Value = 1.251273155212402
$ws.Range("C1").Value = 1.593393325805664
$ws.Range("D1").Value = 3.17042350769043
$ws.Range("E1").Value = -1
